# Auto-generated edit script applying the cryptos.xlsx price/volume update
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '63.061.29'
$ws.Range("E2").Value = '  -1.09%  '
$ws.Range("D3").Value = '2.551.25'
$ws.Range("E3").Value = '  -0.27%  '
$ws.Range("E4").Value = '  -0.02%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '578.04'
$ws.Range("E5").Value = '  +0.64%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '147.01'
$ws.Range("E6").Value = '  -2.08%  '
$ws.Range("E7").Value = '  -0.02%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.583'
$ws.Range("E8").Value = '  -0.68%  '
$ws.Range("E9").Value = '  -1.31%  '
$ws.Range("E10").Value = '  -4.75%  '
$ws.Range("E11").Value = '  -0.63%  '
$ws.Range("E12").Value = '  -1.03%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '27.18'
$ws.Range("E13").Value = '  -3.76%  '
$ws.Range("D14").Value = '3.005.58'
$ws.Range("E14").Value = '  -0.30%  '
$ws.Range("D15").Value = '62.971.78'
$ws.Range("E15").Value = '  -1.07%  '
$ws.Range("E16").Value = '  -0.95%  '
$ws.Range("D17").Value = '2.550.04'
$ws.Range("E17").Value = '  -0.15%  '
$ws.Range("E18").Value = '  -2.15%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '335.82'
$ws.Range("E19").Value = '  -2.18%  '
$ws.Range("E20").Value = '  -0.96%  '
$ws.Range("E21").Value = '  -2.33%  '
$ws.Range("E22").Value = '  +0.03%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '65.37'
$ws.Range("E23").Value = '  -1.21%  '
$ws.Range("E24").Value = '  -0.57%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '1.60'
$ws.Range("E25").Value = '  +1.00%  '
$ws.Range("E26").Value = '  +0.01%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '1.49'
$ws.Range("E27").Value = '  +4.02%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '7.36'
$ws.Range("E29").Value = '  +3.00%  '
$ws.Range("E30").Value = '  +0.78%  '
$ws.Range("D31").Value = '0.0₃0813'
$ws.Range("E31").Value = '  -3.36%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '178.10'
$ws.Range("E32").Value = '  +0.42%  '
$ws.Range("E33").Value = '  -3.98%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '406.71'
$ws.Range("E34").Value = '  -3.91%  '
$ws.Range("B35").Value = 'EthereumClassic'
$ws.Range("C35").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '19.12'
$ws.Range("E35").Value = '  -0.17%  '
$ws.Range("B36").Value = 'PolygonEcosystemToken'
$ws.Range("C36").Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.400'
$ws.Range("E36").Value = '  -2.03%  '
$ws.Range("E37").Value = '  +0.04%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '4.34'
$ws.Range("E38").Value = '  -3.40%  '
$ws.Range("E39").Value = '  -1.27%  '
$ws.Range("E40").Value = '  -0.04%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '39.52'
$ws.Range("E41").Value = '  -2.43%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '151.23'
$ws.Range("E42").Value = '  -3.45%  '
$ws.Range("E43").Value = '  -1.48%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '20.88'
$ws.Range("E44").Value = '  -0.96%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0542'
$ws.Range("E45").Value = '  +1.39%  '
$ws.Range("E46").Value = '  -1.70%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0965'
$ws.Range("E48").Value = '  +2.07%  '
$ws.Range("E49").Value = '  -3.39%  '
$ws.Range("E50").Value = '  -9.19%  '
$ws.Range("E51").Value = '  +0.33%  '
